# Update the fitting parameter "h_p_star" (column K, row 2) on the
# "Parameters" sheet: 0.2533... -> 0.2655...
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Range("K2").Value = 0.26550000000000001
